# Update the division-practice worksheet table with a new set of problems.
# Each populated row in the 5-column table gets new text per cell, while
# preserving the existing run formatting (font, size, etc).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row => list of new cell values (left to right)
$updates = @{
    1  = @("62÷2=", "85÷9=", "74÷9=", "67÷2=", "95÷5=")
    5  = @("50÷4=", "42÷7=", "20÷9=", "34÷2=", "77÷5=")
    9  = @("55÷5=", "36÷3=", "72÷5=", "18÷8=", "50÷3=")
    13 = @("98÷6=", "83÷9=", "10÷9=", "35÷8=", "94÷6=")
    17 = @("97÷3=", "97÷8=", "76÷9=", "34÷2=", "93÷9=")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
